# Rename all 50 worksheets to new "summ" identifiers, preserving order,
# sheetId and rId mapping (only the <sheet name="..."> values change).
$wb = $excel.ActiveWorkbook

$newNames = @(
    "summ53141604",
    "summ53563464",
    "summ54014588",
    "summ54525936",
    "summ54972909",
    "summ55404366",
    "summ55809240",
    "summ56205807",
    "summ56656559",
    "summ57051091",
    "summ57414197",
    "summ57825008",
    "summ58182153",
    "summ58592919",
    "summ59054150",
    "summ59480774",
    "summ59900746",
    "summ00365448",
    "summ00782109",
    "summ01196649",
    "summ01592178",
    "summ01975913",
    "summ02380187",
    "summ02782396",
    "summ03238978",
    "summ03645853",
    "summ04029247",
    "summ04446093",
    "summ04879340",
    "summ05295720",
    "summ05712964",
    "summ06143505",
    "summ06581300",
    "summ06964649",
    "summ07359194",
    "summ07837785",
    "summ08261804",
    "summ08684203",
    "summ09079929",
    "summ09464994",
    "summ09910397",
    "summ10334497",
    "summ10707561",
    "summ11077133",
    "summ11498609",
    "summ11898023",
    "summ12276665",
    "summ12703607",
    "summ13121114",
    "summ13562147"
)

for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $ws.Name = $newNames[$i - 1]
}

